# Auto commit at 2025-08-23 7:58:32.16
#
# Updates the monthly "Metrics" figures (charging income / service income /
# electricity / orders, etc.) with the latest numbers and refreshes the
# dependent "today" comparison sheet. A new helper column (I, and J on the
# electricity row) is staged on the "today" sheet for upcoming work.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Metrics sheet - refresh the raw figures in column B (rows 2-13).
#    Every formula on the "today" sheet references these cells, so the
#    comparison table recalculates automatically once these are written.
# ---------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("Metrics")

$metricUpdates = @{
    "B2"  = 370564.47
    "B3"  = 317929.64
    "B4"  = 117116.20000000001
    "B5"  = 14594
    "B6"  = 3766193.04
    "B7"  = 3196644.3000000003
    "B8"  = 1080758.76
    "B9"  = 145282
    "B10" = 32231516.839999996
    "B11" = 19226514.370000001
    "B12" = 11362467.65
    "B13" = 1242909
}

foreach ($addr in $metricUpdates.Keys) {
    $wsMetrics.Range($addr).Value = $metricUpdates[$addr]
}

# ---------------------------------------------------------------------
# 2) "today" sheet - stage a new column (I) next to the existing B/E/F
#    comparison columns for rows 11-22 (same style as the existing E/F
#    helper cells), plus an extra staging cell (J) on the electricity
#    row (15). B11:F22 formulas pick up the refreshed Metrics values on
#    their own via recalculation.
# ---------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("today")

# Match the numeric formatting already used by the adjacent E/F helper
# columns so the new cells share the same style entry.
$helperFormat = $wsToday.Range("E11").NumberFormat

11..22 | ForEach-Object {
    $wsToday.Range("I$_").NumberFormat = $helperFormat
}
$wsToday.Range("J15").NumberFormat = $helperFormat

# New column I needs to be wide enough for the staged values (OOXML
# width of 15 characters).
$wsToday.Columns.Item(9).ColumnWidth = 14.33

# ---------------------------------------------------------------------
# 3) Restore the cursor/selection on each sheet (last thing touched in
#    each sheet, so the correct sheet stays the active tab).
# ---------------------------------------------------------------------
$wsToday.Activate()
$wsToday.Range("G18").Select() | Out-Null

$wsMetrics.Activate()
$wsMetrics.Range("E7").Select() | Out-Null
